$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive forecaster used to emit an Excel date-serial in column A (e.g.
# 38352 formatted as "2004-12-31 00:00:00"). The bugfix replaces those
# date values with plain "<year>Q4" text labels - the quarter each revision
# row actually refers to - starting at row 2 (row 1 is the header).
$labels = @(
    "2004Q4", "2005Q4", "2006Q4", "2007Q4", "2008Q4", "2009Q4", "2010Q4",
    "2011Q4", "2012Q4", "2013Q4", "2014Q4", "2015Q4", "2016Q4", "2017Q4",
    "2018Q4", "2019Q4", "2020Q4", "2021Q4", "2022Q4", "2023Q4", "2024Q4"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    # Drop the old "YYYY-MM-DD HH:MM:SS" date number format so the cell goes
    # back to the default/general format before writing the text label -
    # otherwise the stored text would still display using the stale date
    # format mask.
    $cell.NumberFormat = "general"
    $cell.Value = $labels[$i]
}
